# Update Leave Card - add new leave card entries / periods for 2023-2024
# and record a new leave usage row (row 76) with remarks "6/16,23/2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix/extend the PERIOD (month-end) dates in column A ---
$ws.Range("A70").Value = 44957   # 01/31/2023
$ws.Range("A72").Value = 44985   # 02/28/2023
$ws.Range("A73").Value = 45016   # 03/31/2023

$ws.Range("A74").Value = 45046   # 04/30/2023
$ws.Range("A75").Value = 45077   # 05/31/2023
$ws.Range("A76").Value = 45107   # 06/30/2023
$ws.Range("A77").Value = 45138   # 07/31/2023
$ws.Range("A78").Value = 45169   # 08/31/2023
$ws.Range("A79").Value = 45199   # 09/30/2023
$ws.Range("A80").Value = 45230   # 10/31/2023
$ws.Range("A81").Value = 45260   # 11/30/2023
$ws.Range("A82").Value = 45291   # 12/31/2023
$ws.Range("A83").Value = 45322   # 01/31/2024
$ws.Range("A84").Value = 45351   # 02/29/2024
$ws.Range("A85").Value = 45382   # 03/31/2024
$ws.Range("A86").Value = 45412   # 04/30/2024
$ws.Range("A87").Value = 45443   # 05/31/2024

# --- Fill EARNED (column C) for the newly dated periods ---
$ws.Range("C74").Value = 1.25
$ws.Range("C75").Value = 1.25
$ws.Range("C76").Value = 1.25

# --- Record a leave usage entry on row 76 ---
$ws.Range("B76").Value = "SL(2-0-0)"
$ws.Range("H76").Value = 2
$ws.Range("K76").Value = "6/16,23/2023"
